$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last three data rows (11-13); remaining rows shift up, dimension becomes A1:T10
$ws.Rows("11:13").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1860306666666667
$ws.Range("H2").Value = 0.558092
$ws.Range("I2").Value = 0.0235467122458118
$ws.Range("J2").Value = 0.02354671224581179
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 30.58864766666666
$ws.Range("N2").Value = 91.765943
$ws.Range("O2").Value = 0.3925391465174898
$ws.Range("P2").Value = 0.3925391465174898
$ws.Range("Q2").Value = 5.690426517861778
$ws.Range("R2").Value = 51.213838660756
$ws.Range("S2").Value = 0.009243006328263886
$ws.Range("T2").Value = 0.009243006328263885

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1860306666666667
$ws.Range("H3").Value = 0.558092
$ws.Range("I3").Value = 0.0235467122458118
$ws.Range("J3").Value = 0.02354671224581179
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.690535
$ws.Range("N3").Value = 68.071605
$ws.Range("O3").Value = 0.291183949679193
$ws.Range("P3").Value = 0.291183949679193
$ws.Range("Q3").Value = 4.221135353073334
$ws.Range("R3").Value = 37.99021817766
$ws.Range("S3").Value = 0.0068564246736949
$ws.Range("T3").Value = 0.006856424673694897

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1860306666666667
$ws.Range("H4").Value = 0.558092
$ws.Range("I4").Value = 0.0235467122458118
$ws.Range("J4").Value = 0.02354671224581179
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.64590566666666
$ws.Range("N4").Value = 73.93771699999999
$ws.Range("O4").Value = 0.3162769038033173
$ws.Range("P4").Value = 0.3162769038033172
$ws.Range("Q4").Value = 4.584894261773777
$ws.Range("R4").Value = 41.264048355964
$ws.Range("S4").Value = 0.007447281243853011
$ws.Range("T4").Value = 0.007447281243853007

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.267456
$ws.Range("H5").Value = 12.802368
$ws.Range("I5").Value = 0.5401505045064059
$ws.Range("J5").Value = 0.5401505045064058
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 30.58864766666666
$ws.Range("N5").Value = 91.765943
$ws.Range("O5").Value = 0.3925391465174898
$ws.Range("P5").Value = 0.3925391465174898
$ws.Range("Q5").Value = 130.5357080170027
$ws.Range("R5").Value = 1174.821372153024
$ws.Range("S5").Value = 0.2120302180299361
$ws.Range("T5").Value = 0.212030218029936

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.267456
$ws.Range("H6").Value = 12.802368
$ws.Range("I6").Value = 0.5401505045064059
$ws.Range("J6").Value = 0.5401505045064058
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 22.690535
$ws.Range("N6").Value = 68.071605
$ws.Range("O6").Value = 0.291183949679193
$ws.Range("P6").Value = 0.291183949679193
$ws.Range("Q6").Value = 96.83085972896001
$ws.Range("R6").Value = 871.4777375606402
$ws.Range("S6").Value = 0.157283157323384
$ws.Range("T6").Value = 0.157283157323384

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.267456
$ws.Range("H7").Value = 12.802368
$ws.Range("I7").Value = 0.5401505045064059
$ws.Range("J7").Value = 0.5401505045064058
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 24.64590566666666
$ws.Range("N7").Value = 73.93771699999999
$ws.Range("O7").Value = 0.3162769038033173
$ws.Range("P7").Value = 0.3162769038033172
$ws.Range("Q7").Value = 105.1753180126507
$ws.Range("R7").Value = 946.577862113856
$ws.Range("S7").Value = 0.1708371291530859
$ws.Range("T7").Value = 0.1708371291530858

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.447007666666666
$ws.Range("H8").Value = 10.341023
$ws.Range("I8").Value = 0.4363027832477824
$ws.Range("J8").Value = 0.4363027832477824
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 30.58864766666666
$ws.Range("N8").Value = 91.765943
$ws.Range("O8").Value = 0.3925391465174898
$ws.Range("P8").Value = 0.3925391465174898
$ws.Range("Q8").Value = 105.4393030199654
$ws.Range("R8").Value = 948.9537271796889
$ws.Range("S8").Value = 0.1712659221592898
$ws.Range("T8").Value = 0.1712659221592898

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.447007666666666
$ws.Range("H9").Value = 10.341023
$ws.Range("I9").Value = 0.4363027832477824
$ws.Range("J9").Value = 0.4363027832477824
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 22.690535
$ws.Range("N9").Value = 68.071605
$ws.Range("O9").Value = 0.291183949679193
$ws.Range("P9").Value = 0.291183949679193
$ws.Range("Q9").Value = 78.21444810576833
$ws.Range("R9").Value = 703.930032951915
$ws.Range("S9").Value = 0.1270443676821141
$ws.Range("T9").Value = 0.1270443676821141

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.447007666666666
$ws.Range("H10").Value = 10.341023
$ws.Range("I10").Value = 0.4363027832477824
$ws.Range("J10").Value = 0.4363027832477824
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 24.64590566666666
$ws.Range("N10").Value = 73.93771699999999
$ws.Range("O10").Value = 0.3162769038033173
$ws.Range("P10").Value = 0.3162769038033172
$ws.Range("Q10").Value = 84.95462578494343
$ws.Range("R10").Value = 764.5916320644909
$ws.Range("S10").Value = 0.1379924934063785
$ws.Range("T10").Value = 0.1379924934063785

